$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @(-0.06591355322105676, 1.530871157958665, 9.613058837089685, 3.100493321568309, 3.142550392515095, 37)
    3  = @(-0.1153523995099622, 1.586536423933229, 9.239019619989865, 3.039575565764053, 3.080471596460598, 36)
    4  = @(0.2107882169940999, 1.437393044426323, 7.605898777194823, 2.757879398595019, 2.789957561700191, 35)
    5  = @(-0.05894863319698449, 1.449923016958758, 8.054208459716339, 2.83799373849139, 2.88005123914993, 34)
    6  = @(0.06660901963189692, 1.455420382974703, 7.924301366729891, 2.815013564217745, 2.857859404195231, 33)
    7  = @(0.01281360325436448, 1.487479230018397, 7.948061675552662, 2.819230688601531, 2.864311658502185, 32)
    8  = @(0.09900455323890892, 1.416118458025723, 8.065898923416519, 2.84005262687446, 2.885244112592491, 31)
    9  = @(-0.005229847462883516, 1.510955205997302, 8.306051057092679, 2.882022043130947, 2.931286119834196, 30)
    10 = @(0.1051993359688867, 1.480362500590311, 8.558816931283857, 2.925545578397961, 2.975403639477648, 29)
    11 = @(-0.01357939262990701, 1.607689944068207, 8.915413051994902, 2.985868893972892, 3.040628603738353, 28)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
